# Update countries & provincias Spain
# Applies the refreshed COVID-19 stats snapshot (28 Abril 2020, 05:22) to the
# "Pais" worksheet. The underlying data source re-sorted all countries by
# "Casos totales" (column B) descending, so only the rows whose ranking
# actually changed need their values rewritten in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 05:22"

# Helper to rewrite a full data row: Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
function Set-CountryRow($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

Set-CountryRow 30 "Pakistan" 13915 0 3029 10594 111 0 292
Set-CountryRow 31 "Chile"    13813 0 7327 6288  426 0 198
Set-CountryRow 32 "Japon"    13614 0 1899 11330 300 0 385

$ws.Cells.Item(45, 4).Value = 2834
$ws.Cells.Item(45, 5).Value = 4388

Set-CountryRow 132 "Congo"  207 0 19 180 0 0 8
Set-CountryRow 133 "Ruanda" 207 0 93 114 0 0 0

Set-CountryRow 140 "Liberia"    133 9 25 92 0 4 16
Set-CountryRow 141 "Madagascar" 128 0 75 53 1 0 0
